# "[master] finished php tutorial_06"
# Fill in row 6 of the "Review" sheet with the review results for
# Tutorial_05 (index.php), which were left blank before.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Review")

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 44550
$ws.Range("C6").Value = "Open"
$ws.Range("D6").Value = "Others"
$ws.Range("E6").Value = "Unknown"
$ws.Range("F6").Value = "Tutorial_05 index.php"
$ws.Range("I6").Value = "1) is ok but test with more data in sample files "
$ws.Range("O6").Value = "PyaePyaeHan"
$ws.Range("P6").Value = "Done"
$ws.Range("V6").Value = "12/20/2021 12:19PM"
$ws.Range("W6").Value = "WaiLinOo"

$ws.Activate()
$ws.Range("W6").Select()
